$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet from Sheet7 to Sheet4
$ws.Name = "Sheet4"

# Update all changed cell values per the diff
$ws.Range("D1").Value = 45574
$ws.Range("F1").Value = 0.877565868055556
$ws.Range("G17").Value = "Turning 1"
$ws.Range("H17").Value = "Turning 2"
$ws.Range("I17").Value = "Turning 3"
$ws.Range("J17").Value = "Turning 4"
$ws.Range("K17").Value = "Turning 5"
$ws.Range("L17").Value = "Turning 6"
$ws.Range("M17").Value = "Turning 7"
$ws.Range("N17").Value = "Turning 8"
$ws.Range("O17").Value = "Turning 9"
$ws.Range("P17").Value = "Turning 10"
$ws.Range("Q17").Value = "Turning 11"
$ws.Range("R17").Value = "Turning 12"
$ws.Range("S17").Value = "Turning 13"
$ws.Range("T17").Value = "Drilling"
$ws.Range("U17").Value = "Milling"
$ws.Range("V17").Value = "Turning 14"
$ws.Range("W17").Value = "Surface Grinding"
$ws.Range("X17").Value = "Induction Hardening"
$ws.Range("B19").Value = "Drilling"
$ws.Range("B20").Value = "Induction Hardening"
$ws.Range("E20").Value = 34.074702596165
$ws.Range("W20").Value = 0
$ws.Range("X20").Value = 34.074702596165
$ws.Range("B21").Value = "Milling"
$ws.Range("E21").Value = 0.147161586307053
$ws.Range("U21").Value = 0.147161586307053
$ws.Range("X21").Value = 0
$ws.Range("B22").Value = "Primary Production Steel Billet"
$ws.Range("E22").Value = 63.9589323653056
$ws.Range("F22").Value = 63.9589323653056
$ws.Range("U22").Value = 0
$ws.Range("B23").Value = "Surface Grinding"
$ws.Range("E23").Value = 0.00859644008161318
$ws.Range("F23").Value = 0
$ws.Range("W23").Value = 0.00859644008161318
$ws.Range("B24").Value = "Turning 1"
$ws.Range("E24").Value = 5.33933131517552
$ws.Range("G24").Value = 5.33933131517552
$ws.Range("V24").Value = 0
$ws.Range("B25").Value = "Turning 10"
$ws.Range("E25").Value = 0.572368330451972
$ws.Range("P25").Value = 0.572368330451972
$ws.Range("R25").Value = 0
$ws.Range("B26").Value = "Turning 11"
$ws.Range("E26").Value = 5.93257354351395
$ws.Range("G26").Value = 0
$ws.Range("Q26").Value = 5.93257354351395
$ws.Range("B27").Value = "Turning 12"
$ws.Range("E27").Value = 0.0199663371087897
$ws.Range("K27").Value = 0
$ws.Range("R27").Value = 0.0199663371087897
$ws.Range("B28").Value = "Turning 13"
$ws.Range("B29").Value = "Turning 14"
$ws.Range("E29").Value = 0.00325209278665893
$ws.Range("J29").Value = 0
$ws.Range("V29").Value = 0.00325209278665893
$ws.Range("B30").Value = "Turning 2"
$ws.Range("E30").Value = 1.66272697476039
$ws.Range("H30").Value = 1.66272697476039
$ws.Range("M30").Value = 0
$ws.Range("B31").Value = "Turning 3"
$ws.Range("B32").Value = "Turning 4"
$ws.Range("E32").Value = 0.366200773790757
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0.366200773790757
$ws.Range("B33").Value = "Turning 5"
$ws.Range("E33").Value = 18.5500128851674
$ws.Range("K33").Value = 18.5500128851674
$ws.Range("L33").Value = 0
$ws.Range("B34").Value = "Turning 6"
$ws.Range("E34").Value = 36.5278087000123
$ws.Range("L34").Value = 36.5278087000123
$ws.Range("N34").Value = 0
$ws.Range("B35").Value = "Turning 7"
$ws.Range("E35").Value = 0.0238991004787029
$ws.Range("M35").Value = 0.0238991004787029
$ws.Range("O35").Value = 0
$ws.Range("B36").Value = "Turning 8"
$ws.Range("E36").Value = 5.68617079237593
$ws.Range("N36").Value = 5.68617079237593
$ws.Range("P36").Value = 0
$ws.Range("B37").Value = "Turning 9"
$ws.Range("E37").Value = 0.170243275878355
$ws.Range("O37").Value = 0.170243275878355
$ws.Range("Q37").Value = 0
